$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / 1h-volume-change snapshot.
# Most rows only update the Price (D) and Volume(1h) (E) columns in
# place. Ranks 47-50 (rows 48-51) also roll: BabyDogeCoin drops out of
# the top 50, Cronos/EnergySwap/Algorand each move up one rank, and
# Mantle newly enters at the bottom (row 51).

$ws.Range("D2").Value = '27.204.15'
$ws.Range("E2").Value = '  -0.18%  '

$ws.Range("D3").Value = '1.647.38'
$ws.Range("E3").Value = '  -0.54%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = '218.29'
$ws.Range("E5").Value = '  -0.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.510'
$ws.Range("E6").Value = '  +1.45%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = '0.256'
$ws.Range("E8").Value = '  +0.72%  '

$ws.Range("D9").Value = '0.0628'
$ws.Range("E9").Value = '  +0.04%  '

$ws.Range("D10").Value = '20.24'
$ws.Range("E10").Value = '  +2.95%  '

$ws.Range("D11").Value = '0.0848'
$ws.Range("E11").Value = '  -0.02%  '

$ws.Range("D12").Value = '1.877.97'
$ws.Range("E12").Value = '  -0.34%  '

$ws.Range("D13").Value = '1.637.37'
$ws.Range("E13").Value = '  -1.03%  '

$ws.Range("E14").Value = '  -1.81%  '

$ws.Range("E15").Value = '  +0.98%  '

$ws.Range("D16").Value = '67.77'
$ws.Range("E16").Value = '  +2.50%  '

$ws.Range("D17").Value = '27.167.25'
$ws.Range("E17").Value = '  -0.12%  '

$ws.Range("D18").Value = '0.0₃0739'
$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("D19").Value = '220.63'
$ws.Range("E19").Value = '  -0.54%  '

$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("E21").Value = '  -0.29%  '

$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("E23").Value = '  +2.52%  '

$ws.Range("E24").Value = '  -0.61%  '

$ws.Range("D25").Value = '148.51'
$ws.Range("E25").Value = '  +0.64%  '

$ws.Range("E26").Value = '  +0.16%  '

$ws.Range("E27").Value = '  +0.38%  '

$ws.Range("E28").Value = '  +0.44%  '

$ws.Range("D29").Value = '15.82'
$ws.Range("E29").Value = '  -0.70%  '

$ws.Range("D30").Value = '0.0506'
$ws.Range("E30").Value = '  -1.79%  '

$ws.Range("E31").Value = '  -0.60%  '

$ws.Range("E32").Value = '  -0.99%  '

$ws.Range("D33").Value = '3.04'
$ws.Range("E33").Value = '  +0.65%  '

$ws.Range("E34").Value = '  +0.11%  '

$ws.Range("D35").Value = '1.275.87'
$ws.Range("E35").Value = '  +0.70%  '

$ws.Range("D36").Value = '2.45'
$ws.Range("E36").Value = '  +0.25%  '

$ws.Range("D38").Value = '0.541'
$ws.Range("E38").Value = '  +0.23%  '

$ws.Range("D39").Value = '0.845'
$ws.Range("E39").Value = '  +1.92%  '

$ws.Range("E40").Value = '  +0.05%  '

$ws.Range("E41").Value = '  +0.16%  '

$ws.Range("D42").Value = '2.24'
$ws.Range("E42").Value = '  +7.18%  '

$ws.Range("D43").Value = '5.38'
$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("D44").Value = '1.788.72'
$ws.Range("E44").Value = '  -0.27%  '

$ws.Range("D45").Value = '62.89'
$ws.Range("E45").Value = '  +1.35%  '

$ws.Range("D46").Value = '92.42'
$ws.Range("E46").Value = '  -0.29%  '

$ws.Range("E47").Value = '  -1.57%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.0513'
$ws.Range("E48").Value = '  -0.73%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.70'
$ws.Range("E49").Value = '  +0.43%  '

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.0974'
$ws.Range("E50").Value = '  -0.48%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.406'
$ws.Range("E51").Value = '  -0.06%  '
